# Apply crypto price/volume updates for Fri Sep  8 19:12:02 UTC 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Assign a string while preserving the cell as text, even when the
    # string looks like a number (Excel would otherwise coerce it to a
    # Double and mangle things like "25.859.07"). Temporarily force the
    # Text number format, assign, then restore the original style so no
    # stray style/number-format is left behind.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "25.859.07"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.630.95"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.54%  "
Set-TextValue $ws.Range("D5") "214.11"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("E7").Value = "  +0.47%  "
Set-TextValue $ws.Range("D9") "0.0632"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("E10").Value = "  -0.21%  "
Set-TextValue $ws.Range("D11") "0.0793"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "1.856.00"
$ws.Range("E12").Value = "  +0.07%  "
Set-TextValue $ws.Range("D13") "4.24"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "1.636.65"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "0.0₃0755"
$ws.Range("E16").Value = "  -0.10%  "
Set-TextValue $ws.Range("D17") "62.53"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "25.860.45"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("E20").Value = "  -0.72%  "
Set-TextValue $ws.Range("D21") "193.23"
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("E22").Value = "  +0.03%  "
Set-TextValue $ws.Range("D23") "6.26"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  +1.43%  "
Set-TextValue $ws.Range("D25") "143.45"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("E27").Value = "  +2.54%  "
$ws.Range("E28").Value = "  +0.13%  "
Set-TextValue $ws.Range("D29") "15.39"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("E30").Value = "  +0.24%  "
Set-TextValue $ws.Range("D31") "0.0499"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").Value = "1.137.84"
$ws.Range("E37").Value = "  -0.23%  "
Set-TextValue $ws.Range("D38") "0.545"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("E41").Value = "  +0.62%  "
Set-TextValue $ws.Range("D42") "99.22"
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("E43").Value = "  -2.98%  "
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "1.765.46"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D47") "56.18"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D48") "0.0528"
$ws.Range("E48").Value = "  +3.31%  "
Set-TextValue $ws.Range("D49") "1.44"
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D50") "0.415"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D51") "7.64"
$ws.Range("E51").Value = "  +1.39%  "
